{"js": "// Load the body's paragraphs (with text) so we can locate the ones that\n// need editing by their current content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate paragraphs by their distinctive (pre-edit) text content rather\n// than hard-coded indices, so the script is resilient to minor shifts.\nlet emptyBeforeFan = null; // the stray empty paragraph right after \"34. The diagram...\"\nlet paraA = null;          // \"(a) | Complete the conversion ...\"\nlet paraEnergy = null;     // \"nae cron > ... / an energy enemy. energy\"\nlet paraB = null;          // \"(b) . Explain why the fan ...\"\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"Complete the conversion of energy\") !== -1) {\n    paraA = items[i];\n    // The empty paragraph immediately preceding this one is the stray\n    // blank paragraph that must be removed.\n    if (i > 0 && items[i - 1].text === \"\") {\n      emptyBeforeFan = items[i - 1];\n    }\n  } else if (t.indexOf(\"an energy enemy. energy\") !== -1 || t.indexOf(\"nae cron\") !== -1) {\n    paraEnergy = items[i];\n  } else if (t.indexOf(\"Explain why the fan continued to spin\") !== -1) {\n    paraB = items[i];\n  }\n}\n\n// 1) Remove the extra blank paragraph right before \"(a) ...\".\nif (emptyBeforeFan) {\n  emptyBeforeFan.delete();\n}\n\n// 2) Fix up the OCR'd text of the three remaining paragraphs.\nif (paraA) {\n  paraA.insertText(\n    \"(a) Complete the conversion of energy when the switch is turned on. [1]\",\n    \"Replace\"\n  );\n}\n\nif (paraEnergy) {\n  paraEnergy.insertText(\n    \"meee cron > | | \\u00b0 | | .\\u000ban energy energy. energy\",\n    \"Replace\"\n  );\n}\n\nif (paraB) {\n  paraB.insertText(\n    \"(b) | Explain why the fan continued to spin for a while even after the switch\\u000bwas-turned off. a it}\",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Collect paragraph text up front (1-based, like real Word COM) so we can\n# locate the paragraphs to edit by their distinctive OCR'd content instead\n# of relying on brittle fixed indices.\n$count = $d.Paragraphs.Count\n\n$emptyBeforeFanIndex = 0\n$paraAIndex = 0\n$paraEnergyIndex = 0\n$paraBIndex = 0\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n\n    if ($t.Contains(\"Complete the conversion of energy\")) {\n        $paraAIndex = $i\n        if ($i -gt 1) {\n            $prevText = $d.Paragraphs.Item($i - 1).Range.Text\n            if ($prevText -eq \"`r\") {\n                $emptyBeforeFanIndex = $i - 1\n            }\n        }\n    } elseif ($t.Contains(\"an energy enemy. energy\") -or $t.Contains(\"nae cron\")) {\n        $paraEnergyIndex = $i\n    } elseif ($t.Contains(\"Explain why the fan continued to spin\")) {\n        $paraBIndex = $i\n    }\n}\n\n# 1) Remove the stray blank paragraph right before \"(a) ...\".\nif ($emptyBeforeFanIndex -gt 0) {\n    $d.Paragraphs.Item($emptyBeforeFanIndex).Range.Delete()\n}\n\n# Re-resolve indices after the delete shifted everything below it down by one.\nif ($emptyBeforeFanIndex -gt 0 -and $paraAIndex -gt $emptyBeforeFanIndex) { $paraAIndex-- }\nif ($emptyBeforeFanIndex -gt 0 -and $paraEnergyIndex -gt $emptyBeforeFanIndex) { $paraEnergyIndex-- }\nif ($emptyBeforeFanIndex -gt 0 -and $paraBIndex -gt $emptyBeforeFanIndex) { $paraBIndex-- }\n\n# 2) Fix up the OCR'd text of the three remaining paragraphs.\nif ($paraAIndex -gt 0) {\n    $d.Paragraphs.Item($paraAIndex).Range.Text = \"(a) Complete the conversion of energy when the switch is turned on. [1]\"\n}\n\nif ($paraEnergyIndex -gt 0) {\n    $d.Paragraphs.Item($paraEnergyIndex).Range.Text = \"meee cron > | | \u00b0 | | .`van energy energy. energy\"\n}\n\nif ($paraBIndex -gt 0) {\n    $d.Paragraphs.Item($paraBIndex).Range.Text = \"(b) | Explain why the fan continued to spin for a while even after the switch`vwas-turned off. a it}\"\n}\n"}
